$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A55").Value = 111898660
$ws.Range("B55").Value = 100997
$ws.Range("D55").Value = 'CR'
$ws.Range("E55").Value = 223246
$ws.Range("F55").Value = 'Skogsalm'
$ws.Range("G55").Value = 'Ulmus glabra'
$ws.Range("H55").Value = 'Huds.'
$ws.Range("L55").Value = ""
$ws.Range("Q55").Value = 650054
$ws.Range("R55").Value = 6654018
$ws.Range("AC55").Value = 'Stammens omkrets i brösthöjd: 64 cm'
$ws.Range("AJ55").Value = ""
$ws.Range("AK55").Value = ""
$ws.Range("AM55").Value = ""
$ws.Range("AO55").Value = ""
$ws.Range("A56").Value = 111898889
$ws.Range("B56").Value = 98980
$ws.Range("D56").Value = 'LC'
$ws.Range("E56").Value = 222498
$ws.Range("F56").Value = 'Blåsippa'
$ws.Range("G56").Value = 'Hepatica nobilis'
$ws.Range("H56").Value = 'Schreb.'
$ws.Range("K56").Value = 'fullt utvecklade blad'
$ws.Range("L56").Value = ""
$ws.Range("Q56").Value = 650135
$ws.Range("R56").Value = 6654003
$ws.Range("AH56").Value = 'Ängsbarrskog'
$ws.Range("AI56").Value = 'Ungskog'
$ws.Range("AJ56").Value = ""
$ws.Range("AK56").Value = ""
$ws.Range("AM56").Value = ""
$ws.Range("AO56").Value = ""
$ws.Range("A57").Value = 111898507
$ws.Range("B57").Value = 89993
$ws.Range("D57").Value = 'VU'
$ws.Range("E57").Value = 1209
$ws.Range("F57").Value = 'Rynkskinn'
$ws.Range("G57").Value = 'Phlebia centrifuga'
$ws.Range("H57").Value = 'P.Karst.'
$ws.Range("L57").Value = ""
$ws.Range("Q57").Value = 650087
$ws.Range("R57").Value = 6654015
$ws.Range("AC57").Value = ""
$ws.Range("AJ57").Value = 'gran'
$ws.Range("AK57").Value = 'Picea abies'
$ws.Range("AM57").Value = 'Liggande död trädstam, utan markontakt'
$ws.Range("AO57").Value = 'Horizontal, dead without ground contact # Picea abies'
$ws.Range("A58").Value = 111898191
$ws.Range("B58").Value = 90480
$ws.Range("E58").Value = 4769
$ws.Range("F58").Value = 'Svavelriska'
$ws.Range("G58").Value = 'Lactarius scrobiculatus'
$ws.Range("H58").Value = '(Scop.:Fr.) Fr.'
$ws.Range("I58").Value = "'2"
$ws.Range("J58").Value = 'fruktkroppar'
$ws.Range("K58").Value = ""
$ws.Range("L58").Value = ""
$ws.Range("A59").Value = 111898336
$ws.Range("B59").Value = 89553
$ws.Range("D59").Value = 'NT'
$ws.Range("E59").Value = 1202
$ws.Range("F59").Value = 'Ullticka'
$ws.Range("G59").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H59").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I59").Value = ""
$ws.Range("J59").Value = ""
$ws.Range("Q59").Value = 650105
$ws.Range("R59").Value = 6654011
$ws.Range("AH59").Value = 'Ängsblandskog'
$ws.Range("AI59").Value = ""
$ws.Range("AJ59").Value = 'gran'
$ws.Range("AK59").Value = 'Picea abies'
$ws.Range("AM59").Value = 'Liggande död trädstam, utan markontakt'
$ws.Range("AO59").Value = 'Horizontal, dead without ground contact # Picea abies'
$ws.Range("A60").Value = 111911660
$ws.Range("B60").Value = 96735
$ws.Range("I60").Value = "'19"
$ws.Range("Q60").Value = 650027
$ws.Range("R60").Value = 6654299
$ws.Range("A61").Value = 111911698
$ws.Range("B61").Value = 96735
$ws.Range("I61").Value = "'16"
$ws.Range("Q61").Value = 650033
$ws.Range("R61").Value = 6654279
